$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 3
    4  = 2
    5  = -1
    6  = -1
    7  = -3
    8  = 1
    9  = -3
    10 = 3
    12 = 3
    14 = 1
    15 = -1
    17 = 2
    18 = 1
    19 = -3
    20 = 3
    21 = 1
    23 = 3
    24 = 3
    25 = 1
    26 = 7
    28 = -2
    29 = -2
    30 = 3
    32 = 3
    33 = -3
    34 = -1
    35 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
